# Applies the "handles float input" marksheet recalculation edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Summary block (rows 10-12): recalculated Right/Wrong/NotAttempt/Max and
#    Marking/Total numbers after fixing float-handling in the scoring code.
# ---------------------------------------------------------------------------

# Give A10/A11/A12 the same "mtitleStyle" formatting already used by A9,
# by copying its format (keeps the existing style table tidy instead of
# creating duplicate xf records).
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "62/112"

# ---------------------------------------------------------------------------
# 2. Drop the third answer block (columns G:H) completely, and the second
#    answer block (columns D:E) past row 18 - only 3 of its questions are
#    still represented.
# ---------------------------------------------------------------------------
$ws.Range("G15:H40").Clear() | Out-Null
$ws.Range("D19:E40").Clear() | Out-Null

# Fill in the 3 remaining "Student Ans" values for the second block - they
# are all correct, so they pick up the green "correctStyle" formatting
# (copied from B10 which already uses it).
$ws.Range("B10").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------------
# 3. Rebuild the first "Student Ans" column (A16:A40), colour-coded green
#    for correct, red for incorrect, and left blank (normalStyle, untouched)
#    for not-attempted.
# ---------------------------------------------------------------------------
$correct = @(16,17,19,21,22,25,27,30,31,32,33,34,38,39)
$incorrect = @(23,24,28,29,35,37)
# 18,20,26,36,40 stay blank/not attempted - no change needed

$ws.Range("B10").Copy() | Out-Null
foreach ($r in $correct) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}

$ws.Range("C10").Copy() | Out-Null
foreach ($r in $incorrect) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

$answers = @{
    16="Option A"; 17="Option D"; 19="Option C"; 21="Option C"; 22="Option D";
    23="Option A"; 24="Option C"; 25="Option A"; 27="Option A"; 28="Option B";
    29="Option B"; 30="Option B"; 31="Option D"; 32="Option C"; 33="Option D";
    34="Option B"; 35="Option A"; 37="Option B"; 38="Option A"; 39="Option D"
}
foreach ($r in $answers.Keys) {
    $ws.Range("A$r").Value = $answers[$r]
}
